$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update some of the existing data values in column B
$ws.Range("B6").Value = 173
$ws.Range("B7").Value = 111
$ws.Range("B8").Value = 33
$ws.Range("B9").Value = 113

# Remove the old B10 value entirely - that row now gets a couple of
# (empty, underline-styled) marker cells instead, mirroring cells like
# C1 / I3 used elsewhere on the sheet.
$ws.Range("B10").ClearContents()
$ws.Range("C10").Font.Underline = 2
$ws.Range("D10").Font.Underline = 2
$ws.Range("D11").Font.Underline = 2

# Move the active selection down to B10
$ws.Range("B10").Select() | Out-Null

# Configure the page for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
